$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111964863
$ws.Range("B2").Value = 89745
$ws.Range("D2").Value = 'VU'
$ws.Range("E2").Value = 2062
$ws.Range("F2").Value = 'Ulltickeporing'
$ws.Range("G2").Value = 'Skeletocutis brevispora'
$ws.Range("H2").Value = 'Niemelä'
$ws.Range("Q2").Value = 734972.3834676194
$ws.Range("R2").Value = 7088252.533270728
$ws.Range("Z2").Value = '16:12'
$ws.Range("AB2").Value = '16:12'
$ws.Range("A3").Value = 111964050
$ws.Range("B3").Value = 90065
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 898
$ws.Range("F3").Value = 'Blackticka'
$ws.Range("G3").Value = 'Steccherinum collabens'
$ws.Range("H3").Value = '(Fr.) Vesterholt'
$ws.Range("Q3").Value = 734893.3330648565
$ws.Range("R3").Value = 7088354.646951701
$ws.Range("Z3").Value = '15:42'
$ws.Range("AB3").Value = '15:42'
$ws.Range("A5").Value = 111965370
$ws.Range("B5").Value = 81248
$ws.Range("E5").Value = 1312
$ws.Range("F5").Value = 'Gammelgransskål'
$ws.Range("G5").Value = 'Pseudographis pinicola'
$ws.Range("H5").Value = '(Nyl.) Rehm'
$ws.Range("M5").Value = $null
$ws.Range("Q5").Value = 734939.7547518623
$ws.Range("R5").Value = 7088232.371273324
$ws.Range("Z5").Value = '16:38'
$ws.Range("AB5").Value = '16:38'
$ws.Range("A6").Value = 111964556
$ws.Range("B6").Value = 56398
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = 'Tretåig hackspett'
$ws.Range("G6").Value = 'Picoides tridactylus'
$ws.Range("H6").Value = '(Linnaeus, 1758)'
$ws.Range("M6").Value = 'färska spår'
$ws.Range("Q6").Value = 734949.4564622594
$ws.Range("R6").Value = 7088268.525185317
$ws.Range("Z6").Value = '16:06'
$ws.Range("AB6").Value = '16:06'
$ws.Range("A8").Value = 111965439
$ws.Range("B8").Value = 56398
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = 'Tretåig hackspett'
$ws.Range("G8").Value = 'Picoides tridactylus'
$ws.Range("H8").Value = '(Linnaeus, 1758)'
$ws.Range("Q8").Value = 734926.7697699566
$ws.Range("R8").Value = 7088234.05367971
$ws.Range("Z8").Value = '16:40'
$ws.Range("AB8").Value = '16:40'
$ws.Range("A9").Value = 111964175
$ws.Range("B9").Value = 89423
$ws.Range("E9").Value = 5432
$ws.Range("F9").Value = 'Granticka'
$ws.Range("G9").Value = 'Porodaedalea chrysoloma'
$ws.Range("H9").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q9").Value = 734896.4627943118
$ws.Range("R9").Value = 7088342.483217424
$ws.Range("Z9").Value = '15:42'
$ws.Range("AB9").Value = '15:42'
$ws.Range("A10").Value = 111964847
$ws.Range("B10").Value = 89405
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = 'Ullticka'
$ws.Range("G10").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H10").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q10").Value = 734972.3834676194
$ws.Range("R10").Value = 7088252.533270728
$ws.Range("Z10").Value = '16:12'
$ws.Range("AB10").Value = '16:12'
$ws.Range("A11").Value = 111965883
$ws.Range("B11").Value = 55611
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 102612
$ws.Range("F11").Value = 'Järpe'
$ws.Range("G11").Value = 'Tetrastes bonasia'
$ws.Range("H11").Value = '(Linnaeus, 1758)'
$ws.Range("M11").Value = 'lockläte, övriga läten'
$ws.Range("Q11").Value = 734846.6442297549
$ws.Range("R11").Value = 7088238.22626837
$ws.Range("Z11").Value = '17:05'
$ws.Range("AB11").Value = '17:05'
